$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "Teacher ID"
$ws.Range("C1").Value = "Subject Name"
$ws.Range("D1").Value = "Number of Classes"

# Row 2 (Mr. Adams): add Subject Name, shift Number of Classes value
$ws.Range("C2").Value = "Math"
$ws.Range("D2").Value = 1

# Row 3 (Ms. Baker): add Subject Name, shift Number of Classes value
$ws.Range("C3").Value = "Science"
$ws.Range("D3").Value = 1
